$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column A previously used the date-only format (style 3); the
# refreshed export now carries the date+time format (style 2) instead.
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new transfer-control record as row 8.
$ws.Range("A8").Value = 45856
$ws.Range("A8").NumberFormat = "YYYY-MM-DD"
$ws.Range("B8").Value = "asd5678"
$ws.Range("C8").Value = "solo"
$ws.Range("D8").Value = "2025-07-18 13:52:17"
$ws.Range("E8").Value = "2025-07-18 13:52:18"
$ws.Range("F8").Value = "2025-07-18 13:52:19"
$ws.Range("G8").Value = "2025-07-18 13:52:20"
$ws.Range("H8").Value = "2025-07-18 13:52:21"
$ws.Range("I8").Value = "2025-07-18 13:52:22"
$ws.Range("J8").Value = "2025-07-18 13:52:22"
$ws.Range("K8").Value = "0:00:01"
$ws.Range("L8").Value = "0:00:01"
$ws.Range("M8").Value = "0:00:05"
# N8 stays blank (no "Entrada CD" value for this record), same as N5/N6/N7.
$ws.Range("O8").Value = "2025-07-18 13:52:24"
$ws.Range("P8").Value = "2025-07-18 13:52:25"
$ws.Range("Q8").Value = "2025-07-18 13:52:26"
$ws.Range("R8").Value = "2025-07-18 13:52:27"
$ws.Range("S8").Value = "0:00:01"
$ws.Range("T8").Value = "0:00:01"
$ws.Range("U8").Value = "0:00:04"
$ws.Range("V8").Value = "0:00:01"
$ws.Range("W8").Value = "2025-07-18 13:52:23"
